# Adds a new "ODI Bowling Extra" worksheet (scraped bowling attributes) and
# clears out leftover placeholder cells in "ODI Batting Extra" that no longer
# carry any value.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "ODI Bowling Extra" sheet as the last sheet in the book.
# ---------------------------------------------------------------------------
$battingExtraWs = $wb.Worksheets.Item("ODI Batting Extra")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$bowlingExtraWs = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$bowlingExtraWs.Name = "ODI Bowling Extra"

# ---------------------------------------------------------------------------
# 2. Header row, styled the same as every other sheet's header row.
# ---------------------------------------------------------------------------
$headers = @("MATCH_CODE", "MAIDEN_OVERS", "PERCENT_WICKETS_OF_ALL")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $bowlingExtraWs.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Copy the header formatting (bold font, borders, centered) from an existing
# header cell so the new sheet matches the look of the rest of the workbook.
$battingExtraWs.Range("A1:C1").Copy()
$bowlingExtraWs.Range("A1:C1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Data rows.
# ---------------------------------------------------------------------------
$data = @(
    @{Row=2; A="4332"; B=""; C=""},
    @{Row=3; A="4338"; B="1"; C="20.00%"},
    @{Row=4; A="4342"; B="1"; C="10.00%"},
    @{Row=5; A="4345"; B="1"; C="40.00%"},
    @{Row=6; A="4350"; B="2"; C="30.00%"},
    @{Row=7; A="4353"; B=""; C=""},
    @{Row=8; A="4398"; B="0"; C=""},
    @{Row=9; A="4399"; B="2"; C="10.00%"},
    @{Row=10; A="4400"; B="0"; C=""},
    @{Row=11; A="4402"; B=""; C=""},
    @{Row=12; A="4406"; B=""; C=""},
    @{Row=13; A="4410"; B="0"; C=""},
    @{Row=14; A="4435"; B=""; C=""},
    @{Row=15; A="4436"; B=""; C=""},
    @{Row=16; A="4437"; B="0"; C="20.00%"},
    @{Row=17; A="4524"; B=""; C=""},
    @{Row=18; A="4526"; B="0"; C="10.00%"},
    @{Row=19; A="4529"; B="0"; C="20.00%"},
    @{Row=20; A="4609"; B="3"; C="60.00%"},
    @{Row=21; A="4613"; B="1"; C="20.00%"}
)

foreach ($entry in $data) {
    $row = $entry.Row

    $cellA = $bowlingExtraWs.Cells.Item($row, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $entry.A

    if ($entry.B -ne "") {
        $cellB = $bowlingExtraWs.Cells.Item($row, 2)
        $cellB.NumberFormat = "@"
        $cellB.Value = $entry.B
    }

    if ($entry.C -ne "") {
        $cellC = $bowlingExtraWs.Cells.Item($row, 3)
        $cellC.NumberFormat = "@"
        $cellC.Value = $entry.C
    }
}

# ---------------------------------------------------------------------------
# 4. Drop the now-empty placeholder cells left over in "ODI Batting Extra".
# ---------------------------------------------------------------------------
$cellsToClear = @(
    "B2", "C2", "D2", "E2",
    "C3", "D3", "E3",
    "E4",
    "C5", "D5", "E5",
    "B6", "C6", "D6", "E6",
    "E7",
    "C8", "D8", "E8",
    "B9", "C9", "D9", "E9",
    "B10", "C10", "D10", "E10",
    "C11", "D11", "E11",
    "B12", "C12", "D12", "E12",
    "B13", "C13", "D13", "E13",
    "C14", "D14", "E14",
    "B15", "C15", "D15", "E15",
    "C16", "D16", "E16",
    "C18", "D18", "E18",
    "B20", "C20", "D20", "E20", "F20",
    "B21", "C21", "D21", "E21", "F21"
)

foreach ($addr in $cellsToClear) {
    $battingExtraWs.Range($addr).ClearContents()
}

$bowlingExtraWs.Activate()
$bowlingExtraWs.Range("A1").Select()
